$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the obsolete "Lead Number" column (column A) -- every remaining
# header shifts one column to the left.
$ws.Columns("A").Delete()

# Insert a new blank column before the old "Actual Date" column (now F)
# to hold the new "Email Id" header.
$ws.Columns("E").Insert()

# New header text for the inserted column.
$ws.Range("E1").Value = "Email Id"

# "Email Id" and "Actual Date" get the refreshed bold/black header style.
$ws.Range("E1:F1").Font.Bold = $true
$ws.Range("E1:F1").Font.Color = 0

# Give column A (now starting with "Name") an explicit width.
$ws.Columns("A").ColumnWidth = 7.92

# Move the active selection.
$ws.Range("E7").Select() | Out-Null
